# Auto-generated edit script: refreshes the Price (column D) and
# Volume(1h) (column E) values for the cryptos.xlsx symbol list,
# matching the "Updated symbol list" GitHub Actions commit.
#
# Values are written with a leading single-quote (quote-prefix) so that
# Excel keeps them as literal text (matching the original inlineStr
# cells, e.g. "291.60" and "-3.10%") instead of auto-converting them
# to numbers/percentages and losing the exact displayed precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'291.60"
$ws.Range("E2").Value = "'-3.10%"
$ws.Range("D3").Value = "'30.65"
$ws.Range("E3").Value = "'-6.25%"
$ws.Range("D4").Value = "'4.953"
$ws.Range("E4").Value = "'-0.04%"
$ws.Range("D5").Value = "'0.07209"
$ws.Range("E5").Value = "'-5.90%"
$ws.Range("D6").Value = "'1.823"
$ws.Range("E6").Value = "'-6.13%"
$ws.Range("D7").Value = "'7.685"
$ws.Range("E7").Value = "'-1.87%"
$ws.Range("D8").Value = "'3.763"
$ws.Range("E8").Value = "'-0.98%"
$ws.Range("D9").Value = "'0.9000"
$ws.Range("E9").Value = "'-2.16%"
$ws.Range("D10").Value = "'0.1650"
$ws.Range("E10").Value = "'-5.91%"
$ws.Range("E11").Value = "'-0.10%"
$ws.Range("D12").Value = "'0.07972"
$ws.Range("E12").Value = "'-7.39%"
$ws.Range("D13").Value = "'0.03037"
$ws.Range("E13").Value = "'-4.17%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("D15").Value = "'0.001502"
$ws.Range("E15").Value = "'-1.12%"
$ws.Range("D16").Value = "'0.005692"
$ws.Range("E16").Value = "'-0.76%"
$ws.Range("D18").Value = "'3.468"
$ws.Range("E18").Value = "'0.22%"
$ws.Range("E19").Value = "'-3.36%"
$ws.Range("D20").Value = "'0.3318"
$ws.Range("E20").Value = "'-0.93%"
$ws.Range("E21").Value = "'-2.19%"
$ws.Range("D22").Value = "'4.034"
$ws.Range("E22").Value = "'-5.62%"
$ws.Range("D23").Value = "'0.2388"
$ws.Range("E23").Value = "'19.73%"
$ws.Range("D24").Value = "'0.04493"
$ws.Range("E24").Value = "'-0.49%"
$ws.Range("E25").Value = "'-0.64%"
$ws.Range("E26").Value = "'-9.06%"
$ws.Range("E27").Value = "'-0.10%"
$ws.Range("D39").Value = "'0.01578"
$ws.Range("E39").Value = "'-6.65%"
$ws.Range("D40").Value = "'0.04402"
$ws.Range("E40").Value = "'-6.13%"
$ws.Range("D41").Value = "'0.007257"
$ws.Range("E41").Value = "'-2.79%"
$ws.Range("D42").Value = "'0.009946"
$ws.Range("D43").Value = "'0.1308"
$ws.Range("E43").Value = "'-3.19%"
$ws.Range("D44").Value = "'0.002013"
$ws.Range("E44").Value = "'-13.66%"
$ws.Range("D45").Value = "'0.009509"
$ws.Range("E45").Value = "'-9.65%"
$ws.Range("D46").Value = "'0.00006001"
$ws.Range("E46").Value = "'-4.11%"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'2.247"
$ws.Range("E48").Value = "'173.85%"
$ws.Range("D49").Value = "'0.003001"
$ws.Range("E49").Value = "'-3.33%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.19%"

Write-Host "Updated 68 cells"
